$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, using same style as other header cells (copy from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Fill H2:H8 with 0
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
